# Reserva_salida1.xlsx - "Primera verificacion hasta el punto 10"
#
# 1. Mayor_maxima.prn (3rd sheet): drop the duplicate scenario-206 row (row 4),
#    leaving the scenario-211 row shifted up into row 4 and the used range
#    shrinking from A1:I5 to A1:I4.
# 2. Reserva.err (6th sheet): replace the lone, still-empty header row
#    (A1:D1) with an actual error log - a header ("Error") followed by 17
#    rows of error messages (the first seven messages duplicated), giving a
#    used range of A1:A18.

$wb = $excel.ActiveWorkbook

# --- 1. Mayor_maxima.prn ---------------------------------------------------
$wsMayor = $wb.Worksheets.Item("Mayor_maxima.prn")
$wsMayor.Rows.Item(4).Delete()

# --- 2. Reserva.err ---------------------------------------------------------
$wsErr = $wb.Worksheets.Item("Reserva.err")

# Clear the old 4-column header row entirely before laying out the new,
# single-column error list.
$wsErr.Range("A1:D1").ClearContents()

$errorMessages = @(
    "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98011 SGDEHI0713.8",
    "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98012 SGDEHI0813.8",
    "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98013 SGDEHI0913.8",
    "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98014 SGDEHI1013.8",
    "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98015 SGDEHI1113.8",
    "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98016 SGDEHI1213.8",
    "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98017  SGDEHI1413.8",
    "***** ERROR EN LOS DATOS DE reserva_DEMANDAS ***** EL AREA INDICADA COMO 99 U.T.E NO POSEE SYSTEMA",
    "***** ERROR EN LOS DATOS DE reserva_DEMANDAS ***** EL AREA INDICADA COMO 18 A.N.D.E. NO POSEE SYSTEMA",
    "***** ERROR EN LOS DATOS DE reserva_DEMANDAS ***** EL AREA INDICADA COMO 20 BRASIL NO POSEE SYSTEMA"
)

# Row 1 is the header; the first seven error messages each occupy two
# consecutive rows, the last three only one row each.
$wsErr.Cells.Item(1, 1).Value = "Error"

$row = 2
for ($i = 0; $i -lt 7; $i++) {
    $wsErr.Cells.Item($row, 1).Value = $errorMessages[$i]
    $row++
    $wsErr.Cells.Item($row, 1).Value = $errorMessages[$i]
    $row++
}
for ($i = 7; $i -lt 10; $i++) {
    $wsErr.Cells.Item($row, 1).Value = $errorMessages[$i]
    $row++
}
